$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) for several events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 3837
$wsExhibit.Range("F9").Value = 2217
$wsExhibit.Range("F10").Value = 334
$wsExhibit.Range("F12").Value = 737
$wsExhibit.Range("F15").Value = 2117

# Sheet "本地生活" (sheet3): update "想去人数" (column F)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 816
$wsLocal.Range("F4").Value = 2078
$wsLocal.Range("F5").Value = 312

# Sheet "全部类型" (sheet4): update "想去人数" (column F) - combined view of all sheets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 816
$wsAll.Range("F4").Value = 2078
$wsAll.Range("F5").Value = 312
$wsAll.Range("F17").Value = 3837
$wsAll.Range("F24").Value = 2217
$wsAll.Range("F25").Value = 334
$wsAll.Range("F28").Value = 737
$wsAll.Range("F32").Value = 2117
